# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I2").Value = "sv"
$ws.Range("J2").Value = "Statement-opinion"

# Row 4: %/Uninterpretable -> sd/Statement-non-opinion
$ws.Range("I4").Value = "sd"
$ws.Range("J4").Value = "Statement-non-opinion"

# Row 6: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I6").Value = "sv"
$ws.Range("J6").Value = "Statement-opinion"

# Row 9: sd/Statement-non-opinion -> aa/Agree/Accept
$ws.Range("I9").Value = "aa"
$ws.Range("J9").Value = "Agree/Accept"

# Row 12: aa/Agree/Accept -> %/Uninterpretable
$ws.Range("I12").Value = "%"
$ws.Range("J12").Value = "Uninterpretable"

# Row 29: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I29").Value = "sv"
$ws.Range("J29").Value = "Statement-opinion"
